# Add a new 10th slide ("Two Content" layout, same as slide 9) with
# empty Title / Content placeholders -- a blank stub slide for the
# next lesson branch.

$p = $ppt.ActivePresentation

# Reuse the "Two Content" layout already used by slide 9 so the new
# slide gets the same Title + two half-width content placeholders.
$layout = $p.Slides.Item(9).CustomLayout

$s = $p.Slides.AddSlide($p.Slides.Count + 1, $layout)
